# Append a new history row (row 8) to the "historique" worksheet,
# mirroring an "Entrée" (stock-in) movement for "Tournevis cruciforme".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 8

$ws.Cells.Item($newRow, 1).Value = "2025-05-22 12:36:01"
$ws.Cells.Item($newRow, 2).Value = "Tournevis cruciforme"
$ws.Cells.Item($newRow, 3).Value = "Entrée"
$ws.Cells.Item($newRow, 4).Value = 5
$ws.Cells.Item($newRow, 5).Value = 103
$ws.Cells.Item($newRow, 6).Value = 98
